# This script applies the edit described by the diff:
#  - Row 3 and Row 6 swap their entire content (columns A:Q)
#  - Row 7, Row 8, Row 9 rotate: new row7 = old row9, new row8 = old row7, new row9 = old row8
#    (and, while moving, the author name "Lisa B. Davidson" is corrected to "Lisa Davidson")
#
# All cells in this sheet are plain text (inlineStr), including ones that look like
# dates/numbers (e.g. "2023-03-01", "1", "2023"). Excel's COM layer auto-converts such
# strings into real dates/numbers when assigned via Range.Value, so we temporarily force
# the target ranges to Text format before writing, then restore the original style/format
# so no visible formatting changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-FullRows($rowA, $rowB) {
    $rangeA = $ws.Range("A" + $rowA + ":Q" + $rowA)
    $rangeB = $ws.Range("A" + $rowB + ":Q" + $rowB)

    $valA = $rangeA.Value2
    $valB = $rangeB.Value2

    $styleA = $rangeA.Style
    $styleB = $rangeB.Style

    $rangeA.NumberFormat = "@"
    $rangeB.NumberFormat = "@"

    $rangeA.Value = $valB
    $rangeB.Value = $valA

    $rangeA.Style = $styleA
    $rangeB.Style = $styleB
}

# --- Row 3 <-> Row 6 full swap ---
Swap-FullRows 3 6

# --- Row 7 -> Row 8 -> Row 9 -> Row 7 rotation ---
$r7 = $ws.Range("A7:Q7")
$r8 = $ws.Range("A8:Q8")
$r9 = $ws.Range("A9:Q9")

$val7 = $r7.Value2
$val8 = $r8.Value2
$val9 = $r9.Value2

$style7 = $r7.Style
$style8 = $r8.Style
$style9 = $r9.Style

$r7.NumberFormat = "@"
$r8.NumberFormat = "@"
$r9.NumberFormat = "@"

# new row7 = old row9 ; new row8 = old row7 ; new row9 = old row8
$r7.Value = $val9
$r8.Value = $val7
$r9.Value = $val8

$r7.Style = $style7
$r8.Style = $style8
$r9.Style = $style9

# Fix the author name typo that travels with the content moving into row 9
# ("Lisa B. Davidson" -> "Lisa Davidson")
$ws.Range("A9").Value = "Rocio Zapata Bustos, Dawn K. Coletta, Jean‐Philippe Galons, Lisa Davidson, Paul Langlais, Janet L. Funk, Wayne T. Willis, Lawrence J. Mandarino"
